# Update imputed values in the result data sheet (RandomForest algorithm output).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 9.199399999999995
$ws.Range("B6").Value = 6.186500000000001
$ws.Range("B7").Value = 5.201299999999998
$ws.Range("C7").Value = -13.97489999999999
$ws.Range("C12").Value = -11.09809999999999
$ws.Range("E13").Value = 16.7732
$ws.Range("E14").Value = 17.14270000000001
$ws.Range("C15").Value = -14.6233
$ws.Range("B16").Value = 5.3569
$ws.Range("E16").Value = 16.46939999999999
$ws.Range("E19").Value = 16.52239999999999
$ws.Range("B20").Value = 9.110999999999999
$ws.Range("C20").Value = -11.7018
$ws.Range("C21").Value = -11.93790000000001
$ws.Range("C22").Value = -11.7447
$ws.Range("E22").Value = 16.93470000000002
$ws.Range("C23").Value = -12.22970000000001
$ws.Range("B28").Value = 5.517200000000003
$ws.Range("B29").Value = 4.894799999999998
$ws.Range("C29").Value = -11.14260000000001
$ws.Range("B32").Value = 7.281199999999997
$ws.Range("C34").Value = -11.59520000000001
$ws.Range("E36").Value = 16.26930000000001
$ws.Range("B40").Value = 9.38319999999999
$ws.Range("C42").Value = -12.1875
$ws.Range("C43").Value = -13.31949999999999
$ws.Range("C44").Value = -13.7701
$ws.Range("C45").Value = -13.94949999999999
$ws.Range("B46").Value = 5.819399999999998
$ws.Range("C46").Value = -14.23779999999999
$ws.Range("E46").Value = 16.18589999999999
$ws.Range("C50").Value = -14.08999999999999
$ws.Range("E50").Value = 16.70199999999999
$ws.Range("B51").Value = 6.121400000000006
$ws.Range("C51").Value = -12.00500000000001
$ws.Range("B52").Value = 4.892000000000004
$ws.Range("B57").Value = 5.891100000000002
$ws.Range("B59").Value = 5.378399999999999
$ws.Range("B62").Value = 5.520099999999998
$ws.Range("B66").Value = 5.471600000000002
$ws.Range("C66").Value = -11.3277
$ws.Range("C67").Value = -11.3134
$ws.Range("B73").Value = 8.328199999999999
$ws.Range("B74").Value = 9.316999999999988
$ws.Range("C79").Value = -11.35610000000001
$ws.Range("C84").Value = -13.48239999999999
$ws.Range("B92").Value = 5.644999999999995
$ws.Range("C92").Value = -11.3609
$ws.Range("E95").Value = 18.11230000000002
$ws.Range("C97").Value = -11.45430000000001
$ws.Range("E97").Value = 16.627
$ws.Range("B100").Value = 5.614
